$d = $word.ActiveDocument

# --- Programa (Portuguese) & Programa (English) ---------------------------
# Both paragraphs list 10 items separated by ";" with no other semicolons
# anywhere else in the document, so a single global Find/Replace using the
# literal ";" is safe and turns each separator into ";" + manual line break.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(";", $false, $false, $false, $false, $false, $true, 1, $false, ";^l", 2)

# --- Bibliografia -----------------------------------------------------------
# The three references are concatenated back-to-back as "<year>.<n>. " with
# no space; that exact pattern ("digit-period-digit") only occurs at the two
# boundaries between references, so replace it with a manual line break
# inserted right after the closing period of the previous reference.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.MatchWildcards = $true
$find.Execute("([0-9])\.([0-9])", $true, $false, $true, $false, $false, $true, 1, $false, "\1.^l\2", 2)
